$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.463.94"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.503.06"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.03"
$ws.Range("E5").Value = "  -2.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.20"
$ws.Range("E6").Value = "  -4.26%  "
$ws.Range("D7").Value = "3.501.48"
$ws.Range("E7").Value = "  -2.40%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.99"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.421"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000217"
$ws.Range("E13").Value = "  -3.00%  "
$ws.Range("D14").Value = "4.096.95"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.45"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.512.15"
$ws.Range("E16").Value = "  -1.91%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "67.570.54"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.09"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "446.25"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.04"
$ws.Range("E22").Value = "  -8.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.619"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.29"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.642.62"
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000126"
$ws.Range("E27").Value = "  +8.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("E28").Value = "  -5.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.20"
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.47"
$ws.Range("E30").Value = "  -4.02%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.53"
$ws.Range("E32").Value = "  -5.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.164"
$ws.Range("E33").Value = "  +3.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.64"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").Value = "3.487.66"
$ws.Range("E35").Value = "  -2.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.83"
$ws.Range("E37").Value = "  -4.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.02"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "178.46"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.16"
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0879"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.38"
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.876"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.34"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.96"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.23"
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.54"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.53"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.988"
$ws.Range("E51").Value = "  -2.51%  "
